$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 14.21340333333333
$ws.Range("H2").Value = 42.64021
$ws.Range("I2").Value = 0.07497543485230342
$ws.Range("J2").Value = 0.07497543485230343
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.188906999999999
$ws.Range("N2").Value = 15.566721
$ws.Range("O2").Value = 0.02571200377994867
$ws.Range("P2").Value = 0.02571200377994868
$ws.Range("Q2").Value = 73.75202805015665
$ws.Range("R2").Value = 663.7682524514099
$ws.Range("S2").Value = 0.001927768664325721
$ws.Range("T2").Value = 0.001927768664325721

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 14.21340333333333
$ws.Range("H3").Value = 42.64021
$ws.Range("I3").Value = 0.07497543485230342
$ws.Range("J3").Value = 0.07497543485230343
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 67.633555
$ws.Range("N3").Value = 202.900665
$ws.Range("O3").Value = 0.3351369029761694
$ws.Range("P3").Value = 0.3351369029761695
$ws.Range("Q3").Value = 961.3029960821833
$ws.Range("R3").Value = 8651.72696473965
$ws.Range("S3").Value = 0.02512703503569252
$ws.Range("T3").Value = 0.02512703503569253

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 14.21340333333333
$ws.Range("H4").Value = 42.64021
$ws.Range("I4").Value = 0.07497543485230342
$ws.Range("J4").Value = 0.07497543485230343
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 128.9862746666667
$ws.Range("N4").Value = 386.958824
$ws.Range("O4").Value = 0.6391510932438819
$ws.Range("P4").Value = 0.6391510932438819
$ws.Range("Q4").Value = 1833.333946301449
$ws.Range("R4").Value = 16500.00551671304
$ws.Range("S4").Value = 0.04792063115228517
$ws.Range("T4").Value = 0.04792063115228518

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.428335
$ws.Range("H5").Value = 88.285005
$ws.Range("I5").Value = 0.1552339127976335
$ws.Range("J5").Value = 0.1552339127976336
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.188906999999999
$ws.Range("N5").Value = 15.566721
$ws.Range("O5").Value = 0.02571200377994867
$ws.Range("P5").Value = 0.02571200377994868
$ws.Range("Q5").Value = 152.700893479845
$ws.Range("R5").Value = 1374.308041318605
$ws.Range("S5").Value = 0.003991374952628976
$ws.Range("T5").Value = 0.003991374952628977

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 29.428335
$ws.Range("H6").Value = 88.285005
$ws.Range("I6").Value = 0.1552339127976335
$ws.Range("J6").Value = 0.1552339127976336
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 67.633555
$ws.Range("N6").Value = 202.900665
$ws.Range("O6").Value = 0.3351369029761694
$ws.Range("P6").Value = 0.3351369029761695
$ws.Range("Q6").Value = 1990.342913780925
$ws.Range("R6").Value = 17913.08622402832
$ws.Range("S6").Value = 0.05202461277187166
$ws.Range("T6").Value = 0.05202461277187168

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 29.428335
$ws.Range("H7").Value = 88.285005
$ws.Range("I7").Value = 0.1552339127976335
$ws.Range("J7").Value = 0.1552339127976336
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 128.9862746666667
$ws.Range("N7").Value = 386.958824
$ws.Range("O7").Value = 0.6391510932438819
$ws.Range("P7").Value = 0.6391510932438819
$ws.Range("Q7").Value = 3795.85130129268
$ws.Range("R7").Value = 34162.66171163412
$ws.Range("S7").Value = 0.0992179250731329
$ws.Range("T7").Value = 0.09921792507313291

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 145.9323983333333
$ws.Range("H8").Value = 437.797195
$ws.Range("I8").Value = 0.7697906523500631
$ws.Range("J8").Value = 0.7697906523500631
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.188906999999999
$ws.Range("N8").Value = 15.566721
$ws.Range("O8").Value = 0.02571200377994867
$ws.Range("P8").Value = 0.02571200377994868
$ws.Range("Q8").Value = 757.2296432386216
$ws.Range("R8").Value = 6815.066789147594
$ws.Range("S8").Value = 0.01979286016299397
$ws.Range("T8").Value = 0.01979286016299398

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 145.9323983333333
$ws.Range("H9").Value = 437.797195
$ws.Range("I9").Value = 0.7697906523500631
$ws.Range("J9").Value = 0.7697906523500631
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 67.633555
$ws.Range("N9").Value = 202.900665
$ws.Range("O9").Value = 0.3351369029761694
$ws.Range("P9").Value = 0.3351369029761695
$ws.Range("Q9").Value = 9869.92688895941
$ws.Range("R9").Value = 88829.34200063467
$ws.Range("S9").Value = 0.2579852551686053
$ws.Range("T9").Value = 0.2579852551686053

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 145.9323983333333
$ws.Range("H10").Value = 437.797195
$ws.Range("I10").Value = 0.7697906523500631
$ws.Range("J10").Value = 0.7697906523500631
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 128.9862746666667
$ws.Range("N10").Value = 386.958824
$ws.Range("O10").Value = 0.6391510932438819
$ws.Range("P10").Value = 0.6391510932438819
$ws.Range("Q10").Value = 18823.27641418875
$ws.Range("R10").Value = 169409.4877276987
$ws.Range("S10").Value = 0.4920125370184638
$ws.Range("T10").Value = 0.4920125370184638
